$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Row 3: update "PERIOD TO EXPIRE" and "LAST UPDATE"
$ws.Range("H3").Value = -99
# Leading apostrophe forces literal text entry so the date-like string is
# not auto-converted into a date serial value (matches the source report's
# plain-text "LAST UPDATE" column).
$ws.Range("I3").Value = "'04-Nov-2025"

# Row 4: update "PERIOD TO EXPIRE" and "LAST UPDATE"
$ws.Range("H4").Value = 700
$ws.Range("I4").Value = "'04-Nov-2025"
